$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 76 (Sr.No 74): "Remove Duplicates from sorted list. Leetcode" / date / who
# ---------------------------------------------------------------------------
$ws.Range("B76").Value = 45769
$ws.Range("B76").NumberFormat = "[$-14009]dd/mm/yyyy;@"

$ws.Range("C76").Value = "Remove Duplicates from sorted list. Leetcode"
$c76Bold = $ws.Range("C76").Characters(37, 8)
$c76Bold.Font.Bold = $true

$ws.Range("G76").Value = "solved and submitted"

$ws.Rows.Item(76).RowHeight = 57.6

# ---------------------------------------------------------------------------
# Row 77 (Sr.No 75): "Merge in between linked list. Leetcode" / Input / Output / who
# ---------------------------------------------------------------------------
$ws.Range("B77").Value = 45769
$ws.Range("B77").NumberFormat = "[$-14009]dd/mm/yyyy;@"

$ws.Range("C77").Value = "Merge in between linked list. Leetcode"
$c77Bold = $ws.Range("C77").Characters(31, 8)
$c77Bold.Font.Bold = $true

$ws.Range("D77").Value = "Input: list1 = [10,1,13,6,9,5], a = 3, b = 4, list2 = [1000000,1000001,1000002]"

# E77 needs a brand-new style: Menlo 9, bold FF262626, medium-left border, wrap + vcenter.
# Start from a cell that is already Menlo-based (D45, style s=7) so the font/border
# mutate in as few steps as possible (avoids spraying extra font/border records).
$ws.Range("E77").Value = "Output: [10,1,13,1000000,1000001,1000002,5]"
$ws.Range("D45").Copy()
$ws.Range("E77").PasteSpecial(-4122)
$ws.Range("E77").Font.Name = "Menlo"
$ws.Range("E77").Font.Size = 9

$e77Rest = $ws.Range("E77").Characters(8, 36)
$e77Rest.Font.Bold = $false
$e77Rest.Font.Size = 9
$e77Rest.Font.Name = "Menlo"

$ws.Range("E77").Font.Bold = $true
$ws.Range("E77").Font.Color = 2500134
$ws.Range("E77").WrapText = $true

$ws.Range("G77").Value = "solved and submitted"

$ws.Rows.Item(77).RowHeight = 57.6

# ---------------------------------------------------------------------------
# Row 78 (Sr.No 76): "Return node when cycle begins. Leetcode" / who
# ---------------------------------------------------------------------------
$ws.Range("B78").Value = 45769
$ws.Range("B78").NumberFormat = "[$-14009]dd/mm/yyyy;@"

$ws.Range("C78").Value = "Return node when cycle begins. Leetcode"
$c78Bold = $ws.Range("C78").Characters(32, 8)
$c78Bold.Font.Bold = $true

$ws.Range("G78").Value = "solved and took help for  submition"

$ws.Rows.Item(78).RowHeight = 72

# ---------------------------------------------------------------------------
# View state: selection + scroll position
# ---------------------------------------------------------------------------
$ws.Range("F78").Select()
$excel.ActiveWindow.ScrollRow = 74
